$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ellipse_desc (row 69)
$ws.Range("B69").Value = "· Has a center point.\n\n· Roundness defined by two axis."

# poly_desc (row 71)
$ws.Range("B71").Value = "· Formed by three or more straight lines (sides) connected in a loop, points are plotted on a plane."

# level_intro_6_2 (row 95)
$ws.Range("B95").Value = "We'll be looking at three particular quadrilaterals that have two pairs of opposite sides that are parallel, and equal-length."

# level_intro_7_0 (row 97)
$ws.Range("B97").Value = "On this level, we'll be looking at some more interesting quadrilaterals. These ones have emphasis on opposite sides and angles."

# level_intro_8_1 (row 99)
$ws.Range("B99").Value = "One helpful tip is to look at these categories in a hierarchy. Going from top to bottom, each level shares attributes from everything above."

# Update the selection view state to match the saved workbook view
$ws.Range("B100").Select()
